$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1172.823553800583
$ws.Range("B2").Value = 19.54705923000972
$ws.Range("C2").Value = 23.45332653045654
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = 2.557929528511238
$ws.Range("F2").Value = "2.6 epochs/min"
